# Update automatico via Actualizar 09-13-2020 03-26-14
# Append two new daily rows (11 & 12 Sep 2020) to the "Condicion_Pacientes"
# table on Hoja1, expanding the table/autofilter range from A1:F182 to
# A1:F184, and move the viewport/selection to reflect the newly active area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$lastRow = 182

# Carry the existing formatting of the last data row down into the two new
# rows (date style on column A, centered-number style on B:F) exactly like a
# user extending the table by filling down before typing the new figures.
$ws.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$ws.Range("A183:F183").PasteSpecial(-4122)
$ws.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$ws.Range("A184:F184").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 183 -> 2020-09-11
$ws.Cells.Item(183, 1).Value = 44085
$ws.Cells.Item(183, 2).Value = 739
$ws.Cells.Item(183, 3).Value = 247
$ws.Cells.Item(183, 4).Value = 670
$ws.Cells.Item(183, 5).Value = 124
$ws.Cells.Item(183, 6).Value = 27

# Row 184 -> 2020-09-12
$ws.Cells.Item(184, 1).Value = 44086
$ws.Cells.Item(184, 2).Value = 2131
$ws.Cells.Item(184, 3).Value = 1087
$ws.Cells.Item(184, 4).Value = 643
$ws.Cells.Item(184, 5).Value = 106
$ws.Cells.Item(184, 6).Value = 31

# Grow the table (and its autofilter) so it covers the new rows too.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$lo.Resize($ws.Range("A1:F184"))

# Reflect where the user ended up after entering the new data.
$ws.Activate()
$ws.Range("F187").Select()
